$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12 was blank aside from its counter in column A. Give it the same
# formatting as the row above it (row 11) before filling in its values,
# so the new attendance row matches the look of its neighbors.
$ws.Range("B11:I11").Copy()
$ws.Range("B12:I12").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Row 10's meeting date/time moves from "9/5 /4:15" to "9/5 /1:00"
$ws.Range("B10").Value = "9/5 /1:00"

# Row 11's meeting date/time stays "9/8/ 4:15"
$ws.Range("B11").Value = "9/8/ 4:15"

# Start the next entry: row 12 becomes the "9/12 /1:00" meeting, held at
# Google Hangout, with everyone marked as attending (A)
$ws.Range("B12").Value = "9/12 /1:00"
$ws.Range("C12").Value = "Google Hangout"
$ws.Range("D12").Value = "A"
$ws.Range("E12").Value = "A"
$ws.Range("F12").Value = "A"
$ws.Range("G12").Value = "A"
$ws.Range("H12").Value = "A"
$ws.Range("I12").Value = "A"

# Move the active selection down to the newly-started row
$ws.Range("I12").Select()
